$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting existing rows 48-52 down to 49-53.
$ws.Rows.Item(48).Insert()

# Copy the style of the date cell (D) from the row above (row 47) into the
# newly inserted row's date cell so the new row matches the existing date
# formatting used throughout the table.
$ws.Range("D47").Copy()
$ws.Range("D48").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 48 with the new data record.
$ws.Cells.Item(48, 1).Value = 11
$ws.Cells.Item(48, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(48, 3).Value = "Bíobío"
$ws.Cells.Item(48, 4).Value = 44610
$ws.Cells.Item(48, 5).Value = 8
$ws.Cells.Item(48, 6).Value = "Fruta"
$ws.Cells.Item(48, 7).Value = 100103
$ws.Cells.Item(48, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(48, 9).Value = 100103002
$ws.Cells.Item(48, 10).Value = "Ciruela"
$ws.Cells.Item(48, 11).Value = "Black Amber"
$ws.Cells.Item(48, 12).Value = "Primera"
$ws.Cells.Item(48, 13).Value = 220
$ws.Cells.Item(48, 14).Value = 11000
$ws.Cells.Item(48, 15).Value = 12000
$ws.Cells.Item(48, 16).Value = 11455
$ws.Cells.Item(48, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(48, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(48, 19).Value = 636
$ws.Cells.Item(48, 20).Value = 18
